$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = '21-09-2024'
$ws.Range("B16").Value = '20:15:00'
$ws.Range("C16").Value = 'Herpertz/Bevo HC HS1'
$ws.Range("D16").Value = 'HC Visé BM HS1'
$ws.Range("E16").Value = 616
$ws.Range("F16").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G16").Value = 'De Heuf'
$ws.Range("H16").Value = 'zaal 1'

# Row 17
$ws.Range("A17").Value = '21-09-2024'
$ws.Range("B17").Value = '20:15:00'
$ws.Range("C17").Value = 'KTSV Eupen (B) HS1'
$ws.Range("D17").Value = 'Biobest/ Sasja HC HS1'
$ws.Range("E17").Value = 584
$ws.Range("F17").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G17").Value = 'Sporthalle Stockbergerweg'
$ws.Range("H17").Value = 'zaal'

# Row 38
$ws.Range("A38").Value = '15-10-2024'
$ws.Range("B38").Value = '20:30:00'
$ws.Range("C38").Value = 'Sezoens Achilles Bocholt HS1'
$ws.Range("D38").Value = 'Biobest/ Sasja HC HS1'
$ws.Range("E38").Value = 1397
$ws.Range("F38").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G38").Value = 'De Damburg (B)'
$ws.Range("H38").Value = 'zaal 1'

# Row 39
$ws.Range("A39").Value = '16-10-2024'
$ws.Range("B39").Value = '20:15:00'
$ws.Range("C39").Value = 'HUBO Handbal HS1'
$ws.Range("D39").Value = 'JD Techniek/ Hurry-up HS1'
$ws.Range("E39").Value = 1278
$ws.Range("F39").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G39").Value = 'Alverberg (B)'
$ws.Range("H39").Value = 'zaal 1'

# Row 40
$ws.Range("A40").Value = '19-10-2024'
$ws.Range("B40").Value = '19:00:00'
$ws.Range("C40").Value = 'Green Park/Handbal Aalsmeer HS1'
$ws.Range("D40").Value = 'LvanRaak Milieu/Handbal Houten HS1'
$ws.Range("E40").Value = 1248
$ws.Range("F40").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G40").Value = 'Sporthal de Bloemhof'
$ws.Range("H40").Value = 'zaal 1'

# Row 41
$ws.Range("A41").Value = '19-10-2024'
$ws.Range("B41").Value = '20:00:00'
$ws.Range("C41").Value = 'KEMBIT-LIONS/Sittardia HS1'
$ws.Range("D41").Value = 'HC Visé BM HS1'
$ws.Range("E41").Value = 1367
$ws.Range("F41").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G41").Value = 'Stadssporthal'
$ws.Range("H41").Value = 'zaal 1'

# Row 42
$ws.Range("A42").Value = '19-10-2024'
$ws.Range("B42").Value = '20:15:00'
$ws.Range("C42").Value = 'KTSV Eupen (B) HS1'
$ws.Range("D42").Value = 'KRAS/Volendam HS1'
$ws.Range("E42").Value = 1307
$ws.Range("F42").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G42").Value = 'Sporthalle Stockbergerweg'
$ws.Range("H42").Value = 'zaal'

# Row 43
$ws.Range("A43").Value = '19-10-2024'
$ws.Range("B43").Value = '20:15:00'
$ws.Range("C43").Value = 'Herpertz/Bevo HC HS1'
$ws.Range("D43").Value = 'Sporting Pelt HS1'
$ws.Range("E43").Value = 1336
$ws.Range("F43").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G43").Value = 'De Heuf'
$ws.Range("H43").Value = 'zaal 1'

# Row 48
$ws.Range("A48").Value = '26-10-2024'
$ws.Range("B48").Value = '20:15:00'
$ws.Range("C48").Value = 'Sporting Pelt HS1'
$ws.Range("D48").Value = 'HUBO Handbal HS1'
$ws.Range("E48").Value = 1456
$ws.Range("F48").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G48").Value = 'Dommelhof (B)'
$ws.Range("H48").Value = 'zaal 1'

# Row 49
$ws.Range("A49").Value = '26-10-2024'
$ws.Range("B49").Value = '20:15:00'
$ws.Range("C49").Value = 'HC Visé BM HS1'
$ws.Range("D49").Value = 'Sezoens Achilles Bocholt HS1'
$ws.Range("E49").Value = 1519
$ws.Range("F49").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G49").Value = 'Hall Omnisport De Visé (B)'
$ws.Range("H49").Value = 'zaal 1'

# Row 88
$ws.Range("A88").Value = '21-12-2024'
$ws.Range("B88").Value = '20:15:00'
$ws.Range("C88").Value = 'KTSV Eupen (B) HS1'
$ws.Range("D88").Value = 'KEMBIT-LIONS/Sittardia HS1'
$ws.Range("E88").Value = 2761
$ws.Range("F88").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G88").Value = 'Sporthalle Stockbergerweg'
$ws.Range("H88").Value = 'zaal'

# Row 89
$ws.Range("A89").Value = '21-12-2024'
$ws.Range("B89").Value = '20:15:00'
$ws.Range("C89").Value = 'Herpertz/Bevo HC HS1'
$ws.Range("D89").Value = 'KRAS/Volendam HS1'
$ws.Range("E89").Value = 2792
$ws.Range("F89").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G89").Value = 'De Heuf'
$ws.Range("H89").Value = 'zaal 1'

# Row 107
$ws.Range("A107").Value = '15-02-2025'
$ws.Range("B107").Value = '20:15:00'
$ws.Range("C107").Value = 'Sporting Pelt HS1'
$ws.Range("D107").Value = 'Herpertz/Bevo HC HS1'
$ws.Range("E107").Value = 3271
$ws.Range("F107").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G107").Value = 'Dommelhof (B)'
$ws.Range("H107").Value = 'zaal 1'

# Row 108
$ws.Range("A108").Value = '15-02-2025'
$ws.Range("B108").Value = '20:15:00'
$ws.Range("C108").Value = 'HC Visé BM HS1'
$ws.Range("D108").Value = 'KEMBIT-LIONS/Sittardia HS1'
$ws.Range("E108").Value = 3332
$ws.Range("F108").Value = 'Heren Super Handball League Super Handball League'
$ws.Range("G108").Value = 'Hall Omnisport De Visé (B)'
$ws.Range("H108").Value = 'zaal 1'
